$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price/Volume columns so numeric-looking
# strings (e.g. "1.001", "0.9999") are stored as text, matching the
# original inline-string cells instead of being parsed as numbers.
$valueRange = $ws.Range("D2:E51")
$valueRange.NumberFormat = "@"

$ws.Range("D2").Value = "28.111.86"
$ws.Range("E2").Value = "  +1.55%  "
$ws.Range("D3").Value = "1.790.99"
$ws.Range("E3").Value = "  +1.87%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.54%  "
$ws.Range("D5").Value = "324.16"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("D7").Value = "0.4300"
$ws.Range("E7").Value = "  -2.91%  "
$ws.Range("D8").Value = "0.3632"
$ws.Range("E8").Value = "  -2.47%  "
$ws.Range("D9").Value = "44.73"
$ws.Range("E9").Value = "  -3.04%  "
$ws.Range("D10").Value = "0.07529"
$ws.Range("E10").Value = "  -2.94%  "
$ws.Range("E11").Value = "  -1.14%  "
$ws.Range("D12").Value = "0.9993"
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("D13").Value = "21.74"
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").Value = "6.164"
$ws.Range("E14").Value = "  -0.58%  "
$ws.Range("D15").Value = "7.347"
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("D16").Value = "1.777.25"
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("D17").Value = "91.86"
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("D18").Value = "0.00001067"
$ws.Range("E18").Value = "  -1.34%  "
$ws.Range("D19").Value = "0.06347"
$ws.Range("E19").Value = "  +1.62%  "
$ws.Range("D20").Value = "0.9995"
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("D21").Value = "17.27"
$ws.Range("D22").Value = "5.956"
$ws.Range("E22").Value = "  -3.78%  "
$ws.Range("D23").Value = "28.101.41"
$ws.Range("E23").Value = "  +1.37%  "
$ws.Range("D24").Value = "11.42"
$ws.Range("E24").Value = "  -1.92%  "
$ws.Range("D25").Value = "2.161"
$ws.Range("E25").Value = "  -7.51%  "
$ws.Range("D26").Value = "160.47"
$ws.Range("E26").Value = "  +3.85%  "
$ws.Range("D27").Value = "20.41"
$ws.Range("E27").Value = "  -1.80%  "
$ws.Range("D28").Value = "1.985.64"
$ws.Range("E28").Value = "  +1.37%  "
$ws.Range("D29").Value = "2.186"
$ws.Range("E29").Value = "  -7.42%  "
$ws.Range("D30").Value = "127.11"
$ws.Range("E30").Value = "  -1.54%  "
$ws.Range("D31").Value = "1.167"
$ws.Range("E31").Value = "  -3.55%  "
$ws.Range("D32").Value = "5.707"
$ws.Range("E32").Value = "  -1.23%  "
$ws.Range("D33").Value = "0.09007"
$ws.Range("E33").Value = "  -2.87%  "
$ws.Range("D34").Value = "3.514"
$ws.Range("E34").Value = "  -4.56%  "
$ws.Range("D35").Value = "12.69"
$ws.Range("E35").Value = "  -0.41%  "
$ws.Range("D36").Value = "0.02326"
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("D37").Value = "5.102"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "0.6465"
$ws.Range("E38").Value = "  -0.61%  "
$ws.Range("D39").Value = "0.2115"
$ws.Range("E39").Value = "  -3.44%  "
$ws.Range("D40").Value = "0.06068"
$ws.Range("E40").Value = "  -0.99%  "
$ws.Range("E41").Value = "  -0.87%  "
$ws.Range("D42").Value = "1.413"
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("D43").Value = "0.9997"
$ws.Range("E43").Value = "  -0.38%  "
$ws.Range("D44").Value = "7.886"
$ws.Range("E44").Value = "  -1.46%  "
$ws.Range("D45").Value = "13.64"
$ws.Range("E45").Value = "  -1.80%  "
$ws.Range("D46").Value = "0.5994"
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("D47").Value = "3.705"
$ws.Range("E47").Value = "  -1.38%  "
$ws.Range("D48").Value = "124.48"
$ws.Range("E48").Value = "  -1.11%  "
$ws.Range("D49").Value = "1.992"
$ws.Range("E49").Value = "  -0.44%  "
$ws.Range("D50").Value = "1.153"
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("D51").Value = "0.06956"
$ws.Range("E51").Value = "  +0.80%  "

# Restore the original (General) formatting now that the text values
# are committed, so the cell style index is unchanged from before.
$valueRange.ClearFormats()
